$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1").Value = "Red Status: 6 projects"
$ws.Range("G7").Value = "o3: 18"
